$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Sheet1: new "Shetye et al." row (row 6)
# ---------------------------------------------------------------------------
$ws1.Range("A6").Value = "A deep learning approach to early identification of suggested sexual harassment from videos "
$ws1.Range("B6").Value = "Shetye et al."
$ws1.Range("F6").Value = "Sexual assault`nSexual Harassment`nSexual Violence"
$ws1.Range("I6").Value = "https://drive.google.com/drive/folders/1kRt-MisnnVqurdlDY90XfMHH0l5GbMxK"
$ws1.Hyperlinks.Add($ws1.Range("I6"), "https://drive.google.com/drive/folders/1kRt-MisnnVqurdlDY90XfMHH0l5GbMxK") | Out-Null
$ws1.Range("I2").Copy()
$ws1.Range("I6").PasteSpecial(-4122) | Out-Null
$ws1.Rows.Item(6).RowHeight = 45

Write-Output "sheet1 done"
